$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-27"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 06-27)"

# Update the June 2022 figure (row 7) and the overall Total (row 14)
$ws.Range("I7").Value = 130
$ws.Range("I14").Value = 793
